# Apply the "ValueSet Creation for DocumentReference.type and Rules
# modification" edit to the workbook.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Date property: refreshed generation timestamp.
$meta.Range("B8").Value = "2024-04-04T12:58:54+00:00"

# Contact property: real contact details instead of the placeholder.
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)"

# --- "Include from ..." sheet ------------------------------------------
$inc = $wb.Worksheets.Item("Include from c80-doc-typecode")

# Rename the sheet to reflect the new source system (LOINC instead of
# the HL7 c80-doc-typecodes code system).
$inc.Name = "Include from LOINC"

# Drop the old c80 code display text - the LOINC codes retained don't
# carry these descriptions anymore.
$inc.Range("B2").Value = ""
$inc.Range("B3").Value = ""

# Update the System URI to point at LOINC instead of c80-doc-typecodes.
$inc.Range("B5").Value = "http://loinc.org"
